# Update NATMI ligand-receptor pair results (Plau-Igf2r) following recomputation
# with updated ligand/receptor-expressing cell counts (per Dr Hou advice).
# Columns: E/K = expressing cell counts, G/M = avg expression, H/N = total expression,
# I/J, O/P = specificity scores, Q/R = edge weights, S/T = edge specificity.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 32.935331
$ws.Range("H2").Value = 98.805993
$ws.Range("I2").Value = 0.1836164637112342
$ws.Range("J2").Value = 0.1836164637112342
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.86633166666667
$ws.Range("N2").Value = 56.598995
$ws.Range("O2").Value = 0.1413973975846522
$ws.Range("P2").Value = 0.1413973975846523
$ws.Range("Q2").Value = 621.3688781974483
$ws.Range("R2").Value = 5592.319903777036
$ws.Range("S2").Value = 0.02596289012246525
$ws.Range("T2").Value = 0.02596289012246525
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 32.935331
$ws.Range("H3").Value = 98.805993
$ws.Range("I3").Value = 0.1836164637112342
$ws.Range("J3").Value = 0.1836164637112342
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 81.49602766666666
$ws.Range("N3").Value = 244.488083
$ws.Range("O3").Value = 0.6107878536829223
$ws.Range("P3").Value = 0.6107878536829223
$ws.Range("Q3").Value = 2684.098646386824
$ws.Range("R3").Value = 24156.88781748142
$ws.Range("S3").Value = 0.1121507057710329
$ws.Range("T3").Value = 0.1121507057710329
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 32.935331
$ws.Range("H4").Value = 98.805993
$ws.Range("I4").Value = 0.1836164637112342
$ws.Range("J4").Value = 0.1836164637112342
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.11432133333333
$ws.Range("N4").Value = 30.342964
$ws.Range("O4").Value = 0.07580375136704795
$ws.Range("P4").Value = 0.07580375136704795
$ws.Range("Q4").Value = 333.1185209536947
$ws.Range("R4").Value = 2998.066688583252
$ws.Range("S4").Value = 0.01391881676206298
$ws.Range("T4").Value = 0.01391881676206298
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 32.935331
$ws.Range("H5").Value = 98.805993
$ws.Range("I5").Value = 0.1836164637112342
$ws.Range("J5").Value = 0.1836164637112342
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.95103433333334
$ws.Range("N5").Value = 68.853103
$ws.Range("O5").Value = 0.1720109973653774
$ws.Range("P5").Value = 0.1720109973653775
$ws.Range("Q5").Value = 755.8999125606977
$ws.Range("R5").Value = 6803.099213046279
$ws.Range("S5").Value = 0.03158405105567302
$ws.Range("T5").Value = 0.03158405105567303
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 95.562134
$ws.Range("H6").Value = 286.686402
$ws.Range("I6").Value = 0.5327646808765668
$ws.Range("J6").Value = 0.5327646808765667
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.86633166666667
$ws.Range("N6").Value = 56.598995
$ws.Range("O6").Value = 0.1413973975846522
$ws.Range("P6").Value = 0.1413973975846523
$ws.Range("Q6").Value = 1802.906914818444
$ws.Range("R6").Value = 16226.16223336599
$ws.Range("S6").Value = 0.07533153940096428
$ws.Range("T6").Value = 0.07533153940096428
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 95.562134
$ws.Range("H7").Value = 286.686402
$ws.Range("I7").Value = 0.5327646808765668
$ws.Range("J7").Value = 0.5327646808765667
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 81.49602766666666
$ws.Range("N7").Value = 244.488083
$ws.Range("O7").Value = 0.6107878536829223
$ws.Range("P7").Value = 0.6107878536829223
$ws.Range("Q7").Value = 7787.934316349707
$ws.Range("R7").Value = 70091.40884714737
$ws.Range("S7").Value = 0.3254061959506653
$ws.Range("T7").Value = 0.3254061959506652
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 95.562134
$ws.Range("H8").Value = 286.686402
$ws.Range("I8").Value = 0.5327646808765668
$ws.Range("J8").Value = 0.5327646808765667
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.11432133333333
$ws.Range("N8").Value = 30.342964
$ws.Range("O8").Value = 0.07580375136704795
$ws.Range("P8").Value = 0.07580375136704795
$ws.Range("Q8").Value = 966.5461305750588
$ws.Range("R8").Value = 8698.915175175529
$ws.Range("S8").Value = 0.04038556140631192
$ws.Range("T8").Value = 0.04038556140631191
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 95.562134
$ws.Range("H9").Value = 286.686402
$ws.Range("I9").Value = 0.5327646808765668
$ws.Range("J9").Value = 0.5327646808765667
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.95103433333334
$ws.Range("N9").Value = 68.853103
$ws.Range("O9").Value = 0.1720109973653774
$ws.Range("P9").Value = 0.1720109973653775
$ws.Range("Q9").Value = 2193.249818400601
$ws.Range("R9").Value = 19739.24836560541
$ws.Range("S9").Value = 0.0916413841186253
$ws.Range("T9").Value = 0.09164138411862528
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.06916133333333
$ws.Range("H10").Value = 111.207484
$ws.Range("I10").Value = 0.2066628180165514
$ws.Range("J10").Value = 0.2066628180165514
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.86633166666667
$ws.Range("N10").Value = 56.598995
$ws.Range("O10").Value = 0.1413973975846522
$ws.Range("P10").Value = 0.1413973975846523
$ws.Range("Q10").Value = 699.3590923198423
$ws.Range("R10").Value = 6294.23183087858
$ws.Range("S10").Value = 0.02922158464505095
$ws.Range("T10").Value = 0.02922158464505095
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.06916133333333
$ws.Range("H11").Value = 111.207484
$ws.Range("I11").Value = 0.2066628180165514
$ws.Range("J11").Value = 0.2066628180165514
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 81.49602766666666
$ws.Range("N11").Value = 244.488083
$ws.Range("O11").Value = 0.6107878536829223
$ws.Range("P11").Value = 0.6107878536829223
$ws.Range("Q11").Value = 3020.989397601463
$ws.Range("R11").Value = 27188.90457841317
$ws.Range("S11").Value = 0.1262271390523938
$ws.Range("T11").Value = 0.1262271390523938
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 37.06916133333333
$ws.Range("H12").Value = 111.207484
$ws.Range("I12").Value = 0.2066628180165514
$ws.Range("J12").Value = 0.2066628180165514
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.11432133333333
$ws.Range("N12").Value = 30.342964
$ws.Range("O12").Value = 0.07580375136704795
$ws.Range("P12").Value = 0.07580375136704795
$ws.Range("Q12").Value = 374.9294092825085
$ws.Range("R12").Value = 3374.364683542576
$ws.Range("S12").Value = 0.01566581687374014
$ws.Range("T12").Value = 0.01566581687374014
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 37.06916133333333
$ws.Range("H13").Value = 111.207484
$ws.Range("I13").Value = 0.2066628180165514
$ws.Range("J13").Value = 0.2066628180165514
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 22.95103433333334
$ws.Range("N13").Value = 68.853103
$ws.Range("O13").Value = 0.1720109973653774
$ws.Range("P13").Value = 0.1720109973653775
$ws.Range("Q13").Value = 850.7755944692059
$ws.Range("R13").Value = 7656.980350222852
$ws.Range("S13").Value = 0.0355482774453665
$ws.Range("T13").Value = 0.03554827744536651
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.80362366666667
$ws.Range("H14").Value = 41.410871
$ws.Range("I14").Value = 0.07695603739564764
$ws.Range("J14").Value = 0.07695603739564763
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 18.86633166666667
$ws.Range("N14").Value = 56.598995
$ws.Range("O14").Value = 0.1413973975846522
$ws.Range("P14").Value = 0.1413973975846523
$ws.Range("Q14").Value = 260.4237422971828
$ws.Range("R14").Value = 2343.813680674645
$ws.Range("S14").Value = 0.01088138341617175
$ws.Range("T14").Value = 0.01088138341617175
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.80362366666667
$ws.Range("H15").Value = 41.410871
$ws.Range("I15").Value = 0.07695603739564764
$ws.Range("J15").Value = 0.07695603739564763
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 81.49602766666666
$ws.Range("N15").Value = 244.488083
$ws.Range("O15").Value = 0.6107878536829223
$ws.Range("P15").Value = 0.6107878536829223
$ws.Range("Q15").Value = 1124.940496238921
$ws.Range("R15").Value = 10124.46446615029
$ws.Range("S15").Value = 0.04700381290883033
$ws.Range("T15").Value = 0.04700381290883032
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.80362366666667
$ws.Range("H16").Value = 41.410871
$ws.Range("I16").Value = 0.07695603739564764
$ws.Range("J16").Value = 0.07695603739564763
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.11432133333333
$ws.Range("N16").Value = 30.342964
$ws.Range("O16").Value = 0.07580375136704795
$ws.Range("P16").Value = 0.07580375136704795
$ws.Range("Q16").Value = 139.6142853290716
$ws.Range("R16").Value = 1256.528567961644
$ws.Range("S16").Value = 0.005833556324932918
$ws.Range("T16").Value = 0.005833556324932918
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.80362366666667
$ws.Range("H17").Value = 41.410871
$ws.Range("I17").Value = 0.07695603739564764
$ws.Range("J17").Value = 0.07695603739564763
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.95103433333334
$ws.Range("N17").Value = 68.853103
$ws.Range("O17").Value = 0.1720109973653774
$ws.Range("P17").Value = 0.1720109973653775
$ws.Range("Q17").Value = 316.8074406980792
$ws.Range("R17").Value = 2851.266966282713
$ws.Range("S17").Value = 0.01323728474571264
$ws.Range("T17").Value = 0.01323728474571264
